$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update country case data per the 5-Aug-2020 13:13 refresh.
# Values below are the authoritative new figures; row position in the sheet
# follows automatically because the sheet is kept sorted by "Casos totales" (col B) desc.

# Row 4
$ws.Cells.Item(4, 2).Value = 4918789
$ws.Cells.Item(4, 3).Value = 369
$ws.Cells.Item(4, 5).Value = 2275564
$ws.Cells.Item(4, 7).Value = 36
$ws.Cells.Item(4, 8).Value = 160326

# Row 14
$ws.Cells.Item(14, 2).Value = 317483
$ws.Cells.Item(14, 3).Value = 2697
$ws.Cells.Item(14, 4).Value = 274932
$ws.Cells.Item(14, 5).Value = 24749
$ws.Cells.Item(14, 7).Value = 185
$ws.Cells.Item(14, 8).Value = 17802

# Row 42
$ws.Cells.Item(42, 2).Value = 68376
$ws.Cells.Item(42, 3).Value = 126
$ws.Cells.Item(42, 4).Value = 63425
$ws.Cells.Item(42, 5).Value = 4374
$ws.Cells.Item(42, 7).Value = 3
$ws.Cells.Item(42, 8).Value = 577

# Row 43
$ws.Cells.Item(43, 2).Value = 61606
$ws.Cells.Item(43, 3).Value = 254
$ws.Cells.Item(43, 4).Value = 55385
$ws.Cells.Item(43, 5).Value = 5868
$ws.Cells.Item(43, 7).Value = 2
$ws.Cells.Item(43, 8).Value = 353

# Row 44
$ws.Cells.Item(44, 1).Value = "Rumania"
$ws.Cells.Item(44, 2).Value = 56550
$ws.Cells.Item(44, 3).Value = 1309
$ws.Cells.Item(44, 4).Value = 28584
$ws.Cells.Item(44, 5).Value = 25445
$ws.Cells.Item(44, 7).Value = 41
$ws.Cells.Item(44, 8).Value = 2521

# Row 45
$ws.Cells.Item(45, 1).Value = "Paises Bajos"
$ws.Cells.Item(45, 2).Value = 55955
$ws.Cells.Item(45, 4).Value = 0
$ws.Cells.Item(45, 5).Value = 0
$ws.Cells.Item(45, 8).Value = 6150

# Row 79
$ws.Cells.Item(79, 1).Value = "Estado de Palestina"
$ws.Cells.Item(79, 2).Value = 13065
$ws.Cells.Item(79, 3).Value = 295
$ws.Cells.Item(79, 4).Value = 6618
$ws.Cells.Item(79, 5).Value = 6359
$ws.Cells.Item(79, 7).Value = 2
$ws.Cells.Item(79, 8).Value = 88

# Row 80
$ws.Cells.Item(80, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(80, 2).Value = 12856
$ws.Cells.Item(80, 4).Value = 6592
$ws.Cells.Item(80, 5).Value = 5891
$ws.Cells.Item(80, 7).Value = 0
$ws.Cells.Item(80, 8).Value = 373

# Row 85
$ws.Cells.Item(85, 2).Value = 10538
$ws.Cells.Item(85, 3).Value = 106
$ws.Cells.Item(85, 4).Value = 6988
$ws.Cells.Item(85, 5).Value = 3332
$ws.Cells.Item(85, 7).Value = 4
$ws.Cells.Item(85, 8).Value = 218

# Row 133
$ws.Cells.Item(133, 2).Value = 1926
$ws.Cells.Item(133, 3).Value = 8
$ws.Cells.Item(133, 5).Value = 91

# Row 152
$ws.Cells.Item(152, 1).Value = "Malta"
$ws.Cells.Item(152, 2).Value = 926
$ws.Cells.Item(152, 3).Value = 36
$ws.Cells.Item(152, 4).Value = 668
$ws.Cells.Item(152, 5).Value = 249
$ws.Cells.Item(152, 8).Value = 9

# Row 153
$ws.Cells.Item(153, 1).Value = "Jamaica"
$ws.Cells.Item(153, 2).Value = 920
$ws.Cells.Item(153, 3).Value = 15
$ws.Cells.Item(153, 4).Value = 745
$ws.Cells.Item(153, 5).Value = 163
$ws.Cells.Item(153, 8).Value = 12

# Row 154
$ws.Cells.Item(154, 1).Value = "Siria"
$ws.Cells.Item(154, 2).Value = 892
$ws.Cells.Item(154, 4).Value = 283
$ws.Cells.Item(154, 5).Value = 563
$ws.Cells.Item(154, 8).Value = 46

# Row 161
$ws.Cells.Item(161, 2).Value = 698
$ws.Cells.Item(161, 3).Value = 26
$ws.Cells.Item(161, 5).Value = 312

# Row 177
$ws.Cells.Item(177, 2).Value = 241
$ws.Cells.Item(177, 3).Value = 14
$ws.Cells.Item(177, 5).Value = 49

# Row 217
$ws.Cells.Item(217, 2).Value = 9
$ws.Cells.Item(217, 3).Value = 1
$ws.Cells.Item(217, 5).Value = 1

# Refresh the "last updated" timestamp banner in A1 (above the header row).
$ws.Range("A1").Value = "Datos actualizados a 5 de Agosto de 2020 a las 13:13"
